# Change Excel Field View to Cache, And set default value to FALSE
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell F1 previously named "View" -> rename to "Cache"
$ws.Range("F1").Value = "Cache"

# Column F (rows 2-26) held boolean TRUE for every data row; set them to FALSE
for ($row = 2; $row -le 26; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
}
